# "work on creating link file" - patient-counts.xlsx
#
# 1) Rename the "# sum of calculations" header label to "sum of calculations"
#    (drop the leading "# "). This is cell E2 on the single worksheet; the
#    shared-string table naturally drops the old, now-unused string and
#    appends the new one at the end, which also re-numbers every other
#    shared-string index used by the row-2 header cells.
#
# 2) Resize the data columns: a new (wider) default-ish column width was
#    applied across the sheet, with a handful of columns (B, E, G, H, I)
#    getting their own distinct widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header text -------------------------------------------------------
$ws.Range("E2").Value = "sum of calculations"

# --- 2. Column widths ------------------------------------------------------
# Target "character" column widths (as stored in the workbook XML) together
# with the ColumnWidth value that reproduces them:
#   A: 17.33203125   B: 13.6640625   C: 17.33203125   D: 17.33203125
#   E: 20             F: 17.33203125  G: 14            H: 10.88671875
#   I: 8.5546875
$ws.Columns.Item(1).ColumnWidth = 16.5
$ws.Columns.Item(2).ColumnWidth = 12.833333333333334
$ws.Columns.Item(3).ColumnWidth = 16.5
$ws.Columns.Item(4).ColumnWidth = 16.5
$ws.Columns.Item(5).ColumnWidth = 19.166666666666668
$ws.Columns.Item(6).ColumnWidth = 16.5
$ws.Columns.Item(7).ColumnWidth = 13.166666666666666
$ws.Columns.Item(8).ColumnWidth = 10
$ws.Columns.Item(9).ColumnWidth = 7.666666666666667
